$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared-string value in A1 with the TRUE() boolean formula.
$ws.Range("A1").Formula = "=TRUE()"

# Move the active selection from D8 back to A1.
$ws.Range("A1").Select()
